# This script splits three long run-on paragraphs (the Portuguese and
# English "Programa" text, and the "Bibliografia" text) into multiple
# sentences/items separated by manual line breaks (w:br), matching the
# target diff.  Word's wildcard Find/Replace represents a manual line
# break in the replacement string with "^l", which becomes a plain
# <w:br/> element splitting the run's text into separate <w:t> runs.

$d = $word.ActiveDocument

$replacements = @(
    # Programa - Portuguese
    @{ Find = "industrial.2) Conceitos";              Replace = "industrial.^l2) Conceitos" },
    @{ Find = "P&ID.3) Medição de Pressão";            Replace = "P&ID.^l3) Medição de Pressão" },
    @{ Find = "de pressão.4) Medição de Nível";         Replace = "de pressão.^l4) Medição de Nível" },
    @{ Find = "ultrassônicos.5) Medição de Vazão";      Replace = "ultrassônicos.^l5) Medição de Vazão" },
    @{ Find = "de vazão.6) Medição de Temperatura";     Replace = "de vazão.^l6) Medição de Temperatura" },
    @{ Find = "de temperatura.7) Monitoramento";        Replace = "de temperatura.^l7) Monitoramento" },

    # Programa - English
    @{ Find = "sector.2)Basic Concepts";                Replace = "sector.^l2)Basic Concepts" },
    @{ Find = "diagrams.3)Pressure Measurement";        Replace = "diagrams.^l3)Pressure Measurement" },
    @{ Find = "sensors.4)Level Measurement";            Replace = "sensors.^l4)Level Measurement" },
    @{ Find = "sensors.5)Flow Measurement";             Replace = "sensors.^l5)Flow Measurement" },
    @{ Find = "technologies.6)Temperature Measurement"; Replace = "technologies.^l6)Temperature Measurement" },
    @{ Find = "sensors.7)Bioreactor Monitoring";        Replace = "sensors.^l7)Bioreactor Monitoring" },

    # Bibliografia
    @{ Find = "694 p.Coughanowr";             Replace = "694 p.^lCoughanowr" },
    @{ Find = "0073397894.DORAN";             Replace = "0073397894.^lDORAN" },
    @{ Find = "0122208553.DUNN";              Replace = "0122208553.^lDUNN" },
    @{ Find = "336 p.FRANCHI";                Replace = "336 p.^lFRANCHI" },
    @{ Find = "9788536512174.SCHMIDELL";      Replace = "9788536512174.^lSCHMIDELL" }
)

foreach ($item in $replacements) {
    $range = $d.Content
    $ok = $range.Find.Execute($item.Find, $true, $false, $false, $false, $false, $true, 1, $false, $item.Replace, 2)
    if (-not $ok) {
        Write-Host "NOT FOUND:" $item.Find
    }
}

Write-Host "Done."
